$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts existing rows 16..86 down to 17..87
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C16").Value = 'Metropolitana'
$ws.Range("D16").Value = 44613
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100114007
$ws.Range("G16").Value = 'Jengibre'
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 790
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 16494
$ws.Range("N16").Value = '$/caja 13 kilos'
$ws.Range("O16").Value = 'Perú'
$ws.Range("P16").Value = 1269
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = 'Hortaliza'
